# Hortaliza, Terminal Hortofrutícola Agro Chillán - Cebolla
# Insert two new weekly-report rows (273 and 274) above the existing
# data block that starts at row 273, shifting the rest of the table
# down by two rows (old 273..371 -> new 275..373).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 273.
$ws.Rows("273:274").Insert()

# --- New row 273 ---
$ws.Range("A273").Value = 7
$ws.Range("B273").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C273").Value = "Ñuble"
$ws.Range("D273").Value = 44510
$ws.Range("E273").Value = 16
$ws.Range("F273").Value = 100112004
$ws.Range("G273").Value = "Cebolla"
$ws.Range("H273").Value = "Sin especificar"
$ws.Range("I273").Value = "1a nueva(o)"
$ws.Range("J273").Value = 12000
$ws.Range("K273").Value = 1100
$ws.Range("L273").Value = 1200
$ws.Range("M273").Value = 1150
$ws.Range("N273").Value = "`$/paquete 10 unidades (volumen en unidades)"
$ws.Range("O273").Value = "Región de O'Higgins"
$ws.Range("P273").Value = 115
$ws.Range("Q273").Value = 10
$ws.Range("R273").Value = "Hortaliza"

# --- New row 274 ---
$ws.Range("A274").Value = 7
$ws.Range("B274").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C274").Value = "Ñuble"
$ws.Range("D274").Value = 44510
$ws.Range("E274").Value = 16
$ws.Range("F274").Value = 100112004
$ws.Range("G274").Value = "Cebolla"
$ws.Range("H274").Value = "Sin especificar"
$ws.Range("I274").Value = "1a nueva(o)"
$ws.Range("J274").Value = 12000
$ws.Range("K274").Value = 800
$ws.Range("L274").Value = 900
$ws.Range("M274").Value = 850
$ws.Range("N274").Value = "`$/paquete 10 unidades (volumen en unidades)"
$ws.Range("O274").Value = "Región del Maule"
$ws.Range("P274").Value = 85
$ws.Range("Q274").Value = 10
$ws.Range("R274").Value = "Hortaliza"
